$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("D2").Value = 95
$ws.Range("D3").Value = 450
$ws.Range("D4").Value = 97
$ws.Range("D13").Value = 339
